# feat: add 2022-Q3 data
#
# Plan:
#  1. The existing "2022-Q2" sheet (2nd sheet) is duplicated right after itself;
#     the duplicate keeps the old Q2 fund-holdings data and is renamed back to
#     "2022-Q2" (this becomes the 3rd sheet).
#  2. The original "2022-Q2" sheet (still 2nd in tab order, keeping its original
#     sheetId/rId) is renamed to "2022-Q3" and its data is replaced with the new
#     Q3 fund-holdings table.
#  3. The "总计" (summary) sheet gets its Q2 summary row pushed down to row 3 and
#     a new Q3 summary row written into row 2.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2    = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Duplicate the current "2022-Q2" sheet so its data survives under the same
#    name, then rename the original in place to "2022-Q3".
# ---------------------------------------------------------------------------
$wsQ2.Copy($null, $wsQ2)
$wsQ2.Name = "2022-Q3"

$wsQ2Dup = $wb.Worksheets.Item(3)
$wsQ2Dup.Name = "2022-Q2"

$wsQ3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 2) Overwrite the "2022-Q3" sheet with the new fund-holdings table.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "004450", "嘉实前沿科技沪港深股票",           "12.58", "89.34", "5.35", "0.6730", 5),
    @(1, "011930", "华夏时代前沿一年持有混合A",         "16.60", "89.70", "3.30", "0.5478", 8),
    @(2, "002980", "华夏创新前沿股票",                 "14.74", "89.33", "3.13", "0.4614", 7),
    @(3, "001759", "嘉实成长增强灵活配置混合",           "4.25",  "90.77", "5.12", "0.2176", 5),
    @(4, "011931", "华夏时代前沿一年持有混合C",         "3.52",  "89.70", "3.30", "0.1162", 8),
    @(5, "012447", "华夏互联网龙头混合A",               "1.48",  "89.62", "5.83", "0.0863", 6),
    @(6, "011924", "嘉实港股互联网产业核心资产混合A",     "1.20",  "87.88", "6.62", "0.0794", 3),
    @(7, "012448", "华夏互联网龙头混合C",               "1.13",  "89.62", "5.83", "0.0659", 6),
    @(8, "011925", "嘉实港股互联网产业核心资产混合C",     "0.41",  "87.88", "6.62", "0.0271", 3)
)

$r = 2
foreach ($row in $rows) {
    $wsQ3.Cells.Item($r, 1).Value = $row[0]

    # Fund code keeps leading zeros -> must stay text.
    $wsQ3.Cells.Item($r, 2).NumberFormat = "@"
    $wsQ3.Cells.Item($r, 2).Value = $row[1]

    $wsQ3.Cells.Item($r, 3).Value = $row[2]

    # Scale / position / weight / market-value columns are stored as text in
    # the source data (e.g. "12.58"), not numbers - force text so trailing
    # zeros / formatting are preserved exactly.
    $wsQ3.Cells.Item($r, 4).NumberFormat = "@"
    $wsQ3.Cells.Item($r, 4).Value = $row[3]
    $wsQ3.Cells.Item($r, 5).NumberFormat = "@"
    $wsQ3.Cells.Item($r, 5).Value = $row[4]
    $wsQ3.Cells.Item($r, 6).NumberFormat = "@"
    $wsQ3.Cells.Item($r, 6).Value = $row[5]
    $wsQ3.Cells.Item($r, 7).NumberFormat = "@"
    $wsQ3.Cells.Item($r, 7).Value = $row[6]

    $wsQ3.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# Match the header / index-column look used on the "总计" sheet.
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the old Q2 row down to row 3 and
#    put the new Q3 totals in row 2.
# ---------------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 9
$wsTotal.Cells.Item(2, 4).Value = 2.27

$wsTotal.Cells.Item(3, 1).Value = 1
